$d = $word.ActiveDocument

$replacements = @(
    @{old="72×70=5040"; new="68×45=3060"},
    @{old="55×38=2090"; new="60×72=4320"},
    @{old="67×46=3082"; new="77×83=6391"},
    @{old="69×63=4347"; new="97×96=9312"},
    @{old="64×87=5568"; new="86×79=6794"},
    @{old="29×51=1479"; new="41×87=3567"},
    @{old="53×62=3286"; new="76×38=2888"},
    @{old="72×78=5616"; new="27×68=1836"},
    @{old="94×74=6956"; new="11×23=253"},
    @{old="53×25=1325"; new="36×37=1332"},
    @{old="42×87=3654"; new="56×68=3808"},
    @{old="84×35=2940"; new="68×43=2924"},
    @{old="51×97=4947"; new="41×72=2952"},
    @{old="71×45=3195"; new="30×81=2430"},
    @{old="86×31=2666"; new="20×43=860"},
    @{old="66×38=2508"; new="39×12=468"},
    @{old="69×71=4899"; new="64×34=2176"},
    @{old="82×93=7626"; new="61×33=2013"},
    @{old="23×22=506"; new="58×79=4582"},
    @{old="59×86=5074"; new="47×76=3572"},
    @{old="97×56=5432"; new="20×25=500"},
    @{old="96×68=6528"; new="70×76=5320"},
    @{old="26×23=598"; new="68×16=1088"},
    @{old="42×39=1638"; new="16×29=464"},
    @{old="95×99=9405"; new="16×41=656"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
